$d = $word.ActiveDocument

# Change 1: "that specific attacks." -> "that specific kind of attack."
$d.Content.Find.Execute("that specific attacks.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "that specific kind of attack.", 2)

# Change 2: "include if it had" -> "include whether it had"
$d.Content.Find.Execute("include if it had", $true, $false, $false, $false, $false,
                         $true, 1, $false, "include whether it had", 2)
